$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values in column D are written with a leading apostrophe so Excel
# stores them as literal text instead of auto-converting numeric-looking
# strings (e.g. "109.05") into floating point numbers.

$ws.Range("D2").Value = "'43.816.12"
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = "'2.300.81"
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'109.05"
$ws.Range("E5").Value = '  +11.86%  '

$ws.Range("D6").Value = "'271.81"
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("E7").Value = '  -0.62%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = '  -1.38%  '

$ws.Range("D10").Value = "'47.20"
$ws.Range("E10").Value = '  +4.78%  '

$ws.Range("E11").Value = '  -1.52%  '

$ws.Range("D12").Value = "'8.42"
$ws.Range("E12").Value = '  +5.01%  '

$ws.Range("E13").Value = '  +2.37%  '

$ws.Range("E14").Value = '  +1.85%  '

$ws.Range("D15").Value = "'2.644.19"
$ws.Range("E15").Value = '  -0.31%  '

$ws.Range("D16").Value = "'0.859"
$ws.Range("E16").Value = '  -1.57%  '

$ws.Range("D17").Value = "'2.296.54"
$ws.Range("E17").Value = '  -0.97%  '

$ws.Range("D18").Value = "'43.788.84"
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").Value = "'0.0000111"
$ws.Range("E19").Value = '  +1.48%  '

$ws.Range("D20").Value = "'6.30"
$ws.Range("E20").Value = '  -1.67%  '

$ws.Range("D21").Value = "'72.26"
$ws.Range("E21").Value = '  -1.75%  '

$ws.Range("D22").Value = "'2.48"
$ws.Range("E22").Value = '  +8.39%  '

$ws.Range("D23").Value = "'234.27"
$ws.Range("E23").Value = '  -2.28%  '

$ws.Range("D24").Value = "'2.97"
$ws.Range("E24").Value = '  +16.94%  '

$ws.Range("D25").Value = "'9.22"
$ws.Range("E25").Value = '  -2.18%  '

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").Value = "'11.38"
$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("D28").Value = "'41.17"
$ws.Range("E28").Value = '  +8.29%  '

$ws.Range("D29").Value = "'3.44"
$ws.Range("E29").Value = '  -1.68%  '

$ws.Range("D30").Value = "'2.28"
$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("D31").Value = "'177.99"
$ws.Range("E31").Value = '  +1.36%  '

$ws.Range("D32").Value = "'21.94"
$ws.Range("E32").Value = '  -2.12%  '

$ws.Range("D33").Value = "'0.0915"
$ws.Range("E33").Value = '  +0.29%  '

$ws.Range("D34").Value = "'5.61"
$ws.Range("E34").Value = '  +2.38%  '

$ws.Range("E35").Value = '  +9.05%  '

$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("D37").Value = "'0.114"
$ws.Range("E37").Value = '  +4.30%  '

$ws.Range("D38").Value = "'0.0359"
$ws.Range("E38").Value = '  -1.31%  '

$ws.Range("D39").Value = "'3.73"
$ws.Range("E39").Value = '  +11.28%  '

$ws.Range("E40").Value = '  -2.94%  '

$ws.Range("D41").Value = "'2.34"
$ws.Range("E41").Value = '  -2.23%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = "'1.40"
$ws.Range("E42").Value = '  -1.67%  '

$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = "'67.37"
$ws.Range("E43").Value = '  +7.07%  '

$ws.Range("D44").Value = "'12.17"
$ws.Range("E44").Value = '  -1.55%  '

$ws.Range("E45").Value = '  +3.04%  '

$ws.Range("E46").Value = '  -3.43%  '

$ws.Range("E47").Value = '  -1.67%  '

$ws.Range("E48").Value = '  +2.28%  '

$ws.Range("D49").Value = "'99.50"
$ws.Range("E49").Value = '  -0.81%  '

$ws.Range("D50").Value = "'0.444"
$ws.Range("E50").Value = '  +6.57%  '

$ws.Range("D51").Value = "'1.51"
$ws.Range("E51").Value = '  +8.53%  '
